$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-7 from
# 45184 (2023-09-15) to 45185 (2023-09-16), keeping existing formatting.
foreach ($row in 2..7) {
    $ws.Range("C$row").Value = 45185
}
